# Update automatic: dades i banners [2026-02-05 20:19]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E2').Value = '2026-02-05 20:17:46'
$ws.Range('E3').Value = '2026-02-05 20:17:48'
$ws.Range('E4').Value = '2026-02-05 20:17:51'
$ws.Range('J4').Value = '989.7 hPa'
$ws.Range('O4').Value = '11.5 °C'
$ws.Range('E5').Value = '2026-02-05 20:17:53'
$ws.Range('J5').Value = '989.9 hPa'
$ws.Range('O5').Value = '10.0 °C'
$ws.Range('E6').Value = '2026-02-05 20:17:55'
$ws.Range('H6').Value = '''71%'
$ws.Range('J6').Value = '991.7 hPa'
$ws.Range('E7').Value = '2026-02-05 20:17:58'
$ws.Range('L7').Value = '56.9 km/h - 268º 19:33 TU'
$ws.Range('O7').Value = '10.4 °C'
$ws.Range('E8').Value = '2026-02-05 20:18:00'
$ws.Range('H8').Value = '''85%'
$ws.Range('M8').Value = '14.9 °C 19:39 TU'
$ws.Range('O8').Value = '8.8 °C'
$ws.Range('E9').Value = '2026-02-05 20:18:03'
$ws.Range('O9').Value = '2.3 °C'
$ws.Range('E10').Value = '2026-02-05 20:18:05'
$ws.Range('E11').Value = '2026-02-05 20:18:08'
$ws.Range('J11').Value = '994.7 hPa'
$ws.Range('O11').Value = '0.7 °C'
$ws.Range('E12').Value = '2026-02-05 20:18:10'
$ws.Range('H12').Value = '''87%'
$ws.Range('O12').Value = '10.1 °C'
$ws.Range('E13').Value = '2026-02-05 20:18:13'
$ws.Range('E14').Value = '2026-02-05 20:18:15'
$ws.Range('H14').Value = '''73%'
$ws.Range('I14').Value = '7.4 mm'
$ws.Range('E15').Value = '2026-02-05 20:18:18'
$ws.Range('H15').Value = '''81%'
$ws.Range('J15').Value = '990.4 hPa'
$ws.Range('O15').Value = '8.4 °C'
$ws.Range('E16').Value = '2026-02-05 20:18:20'
$ws.Range('H16').Value = '''98%'
$ws.Range('E17').Value = '2026-02-05 20:18:23'
$ws.Range('M17').Value = '2.6 °C 19:59 TU'
$ws.Range('E18').Value = '2026-02-05 20:18:25'
$ws.Range('E19').Value = '2026-02-05 20:18:28'
$ws.Range('J19').Value = '992.5 hPa'
$ws.Range('E20').Value = '2026-02-05 20:18:31'
$ws.Range('O20').Value = '-1.3 °C'
$ws.Range('E21').Value = '2026-02-05 20:18:33'
$ws.Range('H21').Value = '''82%'
$ws.Range('J21').Value = '990.6 hPa'
$ws.Range('O21').Value = '6.4 °C'
$ws.Range('E22').Value = '2026-02-05 20:18:36'
$ws.Range('O22').Value = '8.8 °C'
$ws.Range('E23').Value = '2026-02-05 20:18:38'
$ws.Range('J23').Value = '989.8 hPa'
$ws.Range('K23').Value = '3.4 MJ/m2'
$ws.Range('E24').Value = '2026-02-05 20:18:40'
$ws.Range('H24').Value = '''76%'
$ws.Range('J24').Value = '988.9 hPa'
$ws.Range('E25').Value = '2026-02-05 20:18:43'
$ws.Range('J25').Value = '993.9 hPa'
$ws.Range('O25').Value = '0.7 °C'
$ws.Range('E26').Value = '2026-02-05 20:18:45'
$ws.Range('O26').Value = '-0.7 °C'
$ws.Range('E27').Value = '2026-02-05 20:18:48'
$ws.Range('J27').Value = '990.1 hPa'
$ws.Range('E28').Value = '2026-02-05 20:18:51'
$ws.Range('J28').Value = '992.9 hPa'
$ws.Range('O28').Value = '2.6 °C'
$ws.Range('E29').Value = '2026-02-05 20:18:53'
$ws.Range('O29').Value = '9.2 °C'
$ws.Range('E30').Value = '2026-02-05 20:18:56'
$ws.Range('E31').Value = '2026-02-05 20:18:58'
$ws.Range('E32').Value = '2026-02-05 20:19:01'
$ws.Range('H32').Value = '''80%'
$ws.Range('E33').Value = '2026-02-05 20:19:03'
$ws.Range('H33').Value = '''86%'
$ws.Range('O33').Value = '9.3 °C'
$ws.Range('E34').Value = '2026-02-05 20:19:06'
$ws.Range('L34').Value = '49.0 km/h - 260º 19:36 TU'
$ws.Range('O34').Value = '4.1 °C'
$ws.Range('E35').Value = '2026-02-05 20:19:08'
$ws.Range('E36').Value = '2026-02-05 20:19:11'
$ws.Range('H36').Value = '''87%'
